$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H64").Value = 3323.75
$ws.Range("I64").Value = 3095
$ws.Range("K64").Value = 3095
$ws.Range("M64").Value = -2847
$ws.Range("H67").Value = 3323.75
$ws.Range("I67").Value = 3095
$ws.Range("K67").Value = 3095
$ws.Range("M67").Value = -2237
$ws.Range("H74").Value = 4666.6665
$ws.Range("I74").Value = 4666.6665
$ws.Range("K74").Value = 4666.6665
$ws.Range("M74").Value = -3730.6665
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 4666.6665
$ws.Range("I77").Value = 4666.6665
$ws.Range("K77").Value = 23333.3325
$ws.Range("M77").Value = -18653.3325
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H86").Value = 967.6667
$ws.Range("I86").Value = 967.6667
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 967.6667
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("M86").Value = 155.3333
$ws.Range("H89").Value = 967.6667
$ws.Range("I89").Value = 967.6667
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 4838.3335
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("M89").Value = 777.6665000000003
$ws.Range("H106").Value = 2199.3333
$ws.Range("I106").Value = 2199.3333
$ws.Range("K106").Value = 2199.3333
$ws.Range("M106").Value = -1568.3333
$ws.Range("H111").Value = 2550
$ws.Range("I111").Value = 1200
$ws.Range("J111").Value = 3900
$ws.Range("K111").Value = 3600
$ws.Range("L111").Value = 11700
$ws.Range("M111").Value = -533
$ws.Range("N111").Value = -17834
$ws.Range("H113").Value = 26719.25
$ws.Range("I113").Value = 30250.572
$ws.Range("K113").Value = 30250.572
$ws.Range("M113").Value = -26996.572
$ws.Range("H138").Value = 2702.491
$ws.Range("I138").Value = 2288.6453
$ws.Range("J138").Value = 3237.0417
$ws.Range("K138").Value = 6865.9359
$ws.Range("L138").Value = 9711.125100000001
$ws.Range("M138").Value = -1725.9359
$ws.Range("N138").Value = -19991.1251

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 5815103
$ws.Range("I2").Value = 11628706
$ws.Range("K2").Value = 11628706
$ws.Range("M2").Value = -11628593
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 500
$ws.Range("K10").Value = 500
$ws.Range("M10").Value = -330
$ws.Range("H32").Value = 3155.72
$ws.Range("I32").Value = 2470.7903
$ws.Range("K32").Value = 2470.7903
$ws.Range("M32").Value = -2183.7903
$ws.Range("H110").Value = 2222.6667
$ws.Range("I110").Value = 1284.4286
$ws.Range("K110").Value = 1284.4286
$ws.Range("M110").Value = 760.5714
$ws.Range("H116").Value = 5815103
$ws.Range("I116").Value = 11628706
$ws.Range("K116").Value = 11628706
$ws.Range("M116").Value = -11626412

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 5815103
$ws.Range("I3").Value = 11628706
$ws.Range("K3").Value = 11628706
$ws.Range("M3").Value = -11628592
$ws.Range("H107").Value = 1416.3846
$ws.Range("I107").Value = 1296.1578
$ws.Range("J107").Value = 1742.7142
$ws.Range("K107").Value = 1296.1578
$ws.Range("L107").Value = 1742.7142
$ws.Range("M107").Value = 623.8422
$ws.Range("N107").Value = -5582.7142
$ws.Range("H134").Value = 7114.48
$ws.Range("I134").Value = 8089.143
$ws.Range("K134").Value = 24267.429
$ws.Range("M134").Value = -21732.429

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 870.625
$ws.Range("I16").Value = 828
$ws.Range("K16").Value = 828
$ws.Range("M16").Value = -541
$ws.Range("H31").Value = 1838.2609
$ws.Range("J31").Value = 2282.5715
$ws.Range("L31").Value = 2282.5715
$ws.Range("N31").Value = -2872.5715
$ws.Range("H34").Value = 1838.2609
$ws.Range("J34").Value = 2282.5715
$ws.Range("L34").Value = 2282.5715
$ws.Range("N34").Value = -2686.5715
$ws.Range("H58").Value = 2072671.1
$ws.Range("I58").Value = 3107283.8
$ws.Range("K58").Value = 3107283.8
$ws.Range("M58").Value = -3107080.8
$ws.Range("H86").Value = 2348.8333
$ws.Range("I86").Value = 1665
$ws.Range("K86").Value = 1665
$ws.Range("M86").Value = -542
$ws.Range("H89").Value = 2348.8333
$ws.Range("I89").Value = 1665
$ws.Range("K89").Value = 8325
$ws.Range("M89").Value = -2709
$ws.Range("H105").Value = 944.1667
$ws.Range("I105").Value = 1022
$ws.Range("J105").Value = 555
$ws.Range("K105").Value = 1022
$ws.Range("L105").Value = 555
$ws.Range("M105").Value = 725
$ws.Range("N105").Value = -4049
$ws.Range("H113").Value = 870.625
$ws.Range("I113").Value = 828
$ws.Range("K113").Value = 828
$ws.Range("M113").Value = 1342
$ws.Range("H122").Value = 2506.5217
$ws.Range("I122").Value = 1777.6111
$ws.Range("K122").Value = 5332.8333
$ws.Range("M122").Value = -2882.8333
$ws.Range("H132").Value = 2519.45
$ws.Range("I132").Value = 1098.091
$ws.Range("J132").Value = 4256.6665
$ws.Range("K132").Value = 3294.273
$ws.Range("L132").Value = 12769.9995
$ws.Range("M132").Value = -764.2729999999997
$ws.Range("N132").Value = -17829.9995
$ws.Range("H134").Value = 1289.8334
$ws.Range("I134").Value = 1298
$ws.Range("K134").Value = 3894
$ws.Range("M134").Value = -1359
$ws.Range("H136").Value = 2072671.1
$ws.Range("I136").Value = 3107283.8
$ws.Range("K136").Value = 9321851.399999999
$ws.Range("M136").Value = -9319301.399999999
$ws.Range("H141").Value = 41744.582
$ws.Range("J141").Value = 64419.285
$ws.Range("L141").Value = 64419.285
$ws.Range("N141").Value = -74779.285

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H12").Value = 333.33334
$ws.Range("J12").Value = 333.33334
$ws.Range("L12").Value = 1000.00002
$ws.Range("N12").Value = -1346.00002
$ws.Range("H131").Value = 774.50507
$ws.Range("J131").Value = 785.0105
$ws.Range("L131").Value = 2355.0315
$ws.Range("N131").Value = -12435.0315
$ws.Range("H132").Value = 1669.8
$ws.Range("J132").Value = 2725
$ws.Range("L132").Value = 24525
$ws.Range("N132").Value = -29585
$ws.Range("H139").Value = 18333.334
$ws.Range("I139").Value = 50500
$ws.Range("K139").Value = 151500
$ws.Range("M139").Value = -146360

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H113").Value = 1399.75
$ws.Range("I113").Value = 1100
$ws.Range("J113").Value = 1499.6666
$ws.Range("K113").Value = 1100
$ws.Range("L113").Value = 1499.6666
$ws.Range("M113").Value = 1070
$ws.Range("N113").Value = -5839.6666

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H46").Value = 2089.3572
$ws.Range("I46").Value = 1393.8
$ws.Range("J46").Value = 2475.7778
$ws.Range("K46").Value = 1393.8
$ws.Range("L46").Value = 2475.7778
$ws.Range("M46").Value = -1205.8
$ws.Range("N46").Value = -2851.7778
$ws.Range("H61").Value = 1953.6666
$ws.Range("I61").Value = 1827.3334
$ws.Range("K61").Value = 1827.3334
$ws.Range("M61").Value = -1625.3334
$ws.Range("H93").Value = 712.7143
$ws.Range("I93").Value = 637.8
$ws.Range("J93").Value = 900
$ws.Range("K93").Value = 637.8
$ws.Range("L93").Value = 900
$ws.Range("M93").Value = 610.2
$ws.Range("N93").Value = -3396
$ws.Range("H113").Value = 1953.6666
$ws.Range("I113").Value = 1827.3334
$ws.Range("K113").Value = 1827.3334
$ws.Range("M113").Value = 342.6666
$ws.Range("H132").Value = 1755.1111
$ws.Range("I132").Value = 1504.7333
$ws.Range("K132").Value = 4514.199900000001
$ws.Range("M132").Value = -1984.199900000001

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630
$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184
